# Generate Report for Handoff
#
# The CI report re-ran: the row for "a8b4ec37-...md" moved from the 2nd
# data row to the last data row (before the static ".localization-config"
# row), the two rows that follow it shift up, and the a8b4ec37 row's
# status/handoff timestamps are refreshed to reflect a brand new handoff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"

$ws.Range("A3").Value = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"

$ws.Range("A4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
    } elseif ($addr -eq "`$A`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-10 01:15:40"
$ws.Range("E2").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$ws.Range("F2").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$ws.Range("G2").Value = "2016-03-10 01:16:20"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-10 01:15:40"
$ws.Range("E3").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$ws.Range("F3").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
$ws.Range("G3").Value = "2016-03-10 01:16:20"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf"
$ws.Range("D4").Value = "2016-03-10 01:19:19"
$ws.Range("E4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws.Range("F4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf"
$ws.Range("G4").Value = "2016-03-10 01:18:28"
$ws.Range("H4").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
    } elseif ($addr -eq "`$C`$2") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
    } elseif ($addr -eq "`$E`$2") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
    } elseif ($addr -eq "`$F`$2") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
    } elseif ($addr -eq "`$C`$3") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
    } elseif ($addr -eq "`$E`$3") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
    } elseif ($addr -eq "`$F`$3") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.zh-cn.xlf"
    } elseif ($addr -eq "`$A`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
    } elseif ($addr -eq "`$C`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf"
    } elseif ($addr -eq "`$E`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
    } elseif ($addr -eq "`$F`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
$ws.Range("D2").Value = "2016-03-10 01:15:46"
$ws.Range("E2").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$ws.Range("F2").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
$ws.Range("G2").Value = "2016-03-10 01:16:38"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
$ws.Range("D3").Value = "2016-03-10 01:15:46"
$ws.Range("E3").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
$ws.Range("F3").Value = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
$ws.Range("G3").Value = "2016-03-10 01:16:38"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf"
$ws.Range("D4").Value = "2016-03-10 01:19:25"
$ws.Range("E4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
$ws.Range("F4").Value = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf"
$ws.Range("G4").Value = "2016-03-10 01:18:44"
$ws.Range("H4").Value = "Include"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq "`$A`$2") {
        $h.TextToDisplay = "ffff5e68522e-1e57-4423-a2e4-b1a65b775524.md"
    } elseif ($addr -eq "`$C`$2") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
    } elseif ($addr -eq "`$E`$2") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
    } elseif ($addr -eq "`$F`$2") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
    } elseif ($addr -eq "`$A`$3") {
        $h.TextToDisplay = "ffffff7a96c373-ffb6-4d51-b7e3-1b863f91784f.md"
    } elseif ($addr -eq "`$C`$3") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
    } elseif ($addr -eq "`$E`$3") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.md"
    } elseif ($addr -eq "`$F`$3") {
        $h.TextToDisplay = "b8a6c97e-a5ff-4bf2-bb1a-51cf7f19d047.a16e12aaccc184848c04ec814ee57caed850f23a.de-de.xlf"
    } elseif ($addr -eq "`$A`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
    } elseif ($addr -eq "`$C`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf"
    } elseif ($addr -eq "`$E`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.md"
    } elseif ($addr -eq "`$F`$4") {
        $h.TextToDisplay = "a8b4ec37-3611-47ba-989c-b0307579875a.b1350e8442248f251d70466f37bcff4c4b5969df.de-de.xlf"
    }
}
